$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 500
$ws.Range("C4").Value = 500
$ws.Range("D4").Formula = "=(C4*100)/B4"
$ws.Range("E4").Value = 0
$ws.Range("F4").Formula = "=(E4*100)/B4"
$ws.Range("G4").Value = 101

$ws.Range("A5").Value = 1
$ws.Range("B5").Value = 100
$ws.Range("C5").Value = 100
$ws.Range("D5").Formula = "=(C5*100)/B5"
$ws.Range("E5").Value = 0
$ws.Range("F5").Formula = "=(E5*100)/B5"
$ws.Range("G5").Value = 101

$ws.Range("A6").Value = 2
$ws.Range("B6").Value = 100
$ws.Range("C6").Value = 59
$ws.Range("D6").Formula = "=(C6*100)/B6"
$ws.Range("E6").Value = 39
$ws.Range("F6").Formula = "=(E6*100)/B6"
$ws.Range("G6").Value = 97.3

$ws.Range("A7").Value = 2
$ws.Range("B7").Value = 500
$ws.Range("C7").Value = 288
$ws.Range("D7").Formula = "=(C7*100)/B7"
$ws.Range("E7").Value = 196
$ws.Range("F7").Formula = "=(E7*100)/B7"
$ws.Range("G7").Value = 97.1

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 500

$ws.Range("C8").Select()
